# Small updates from Eliz's feedback
#
# 1) "...and perhaps catching a variety of stages..." -> "...and perhaps catch a variety of stages..."
# 2) "...value between genomes..." -> "...value, which is a measure of genome differentation, between genomes..."

$d = $word.ActiveDocument

$found1 = $d.Content.Find.Execute(
    "catching a variety of stages of speciation put forth by",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "catch a variety of stages of speciation put forth by", 2
)
Write-Output "replace 1 (perhaps catching -> perhaps catch): $found1"

$found2 = $d.Content.Find.Execute(
    "value between genomes as previously adjusted for metagenomics data",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "value, which is a measure of genome differentation, between genomes as previously adjusted for metagenomics data", 2
)
Write-Output "replace 2 (value between genomes -> value, which is a measure of genome differentation, between genomes): $found2"
